# Apply updated crypto price/volume values to Sheet1, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.232.12"
$ws.Range("E2").Value = "  +8.95%  "
$ws.Range("D3").Value = "3.454.99"
$ws.Range("E3").Value = "  +6.05%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "`'414.30"
$ws.Range("E5").Value = "  +4.10%  "
$ws.Range("D6").Value = "`'124.10"
$ws.Range("E6").Value = "  +14.23%  "
$ws.Range("D7").Value = "3.448.06"
$ws.Range("E7").Value = "  +5.97%  "
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").Value = "`'1.00"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "`'0.659"
$ws.Range("E10").Value = "  +6.49%  "
$ws.Range("D11").Value = "`'0.130"
$ws.Range("E11").Value = "  +35.83%  "
$ws.Range("D12").Value = "`'41.39"
$ws.Range("E12").Value = "  +5.24%  "
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "3.994.47"
$ws.Range("E14").Value = "  +5.87%  "
$ws.Range("D15").Value = "`'8.51"
$ws.Range("E15").Value = "  +3.06%  "
$ws.Range("D16").Value = "`'19.76"
$ws.Range("E16").Value = "  +4.19%  "
$ws.Range("D17").Value = "3.445.70"
$ws.Range("E17").Value = "  +5.90%  "
$ws.Range("D18").Value = "62.171.65"
$ws.Range("E18").Value = "  +9.21%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "`'11.28"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("E21").Value = "  +23.28%  "
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").Value = "`'82.22"
$ws.Range("E23").Value = "  +11.01%  "
$ws.Range("D24").Value = "`'314.49"
$ws.Range("E24").Value = "  +7.17%  "
$ws.Range("D25").Value = "`'13.06"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").Value = "`'3.19"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("D27").Value = "`'31.15"
$ws.Range("E27").Value = "  +11.10%  "
$ws.Range("D28").Value = "`'7.83"
$ws.Range("E28").Value = "  +4.69%  "
$ws.Range("D29").Value = "`'7.90"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "`'4.30"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("D31").Value = "`'0.173"
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("E32").Value = "  +3.94%  "
$ws.Range("D33").Value = "`'11.57"
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "`'42.53"
$ws.Range("E34").Value = "  +6.20%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").Value = "`'2.58"
$ws.Range("E35").Value = "  +20.54%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "`'52.35"
$ws.Range("E38").Value = "  +2.09%  "
$ws.Range("D39").Value = "`'3.51"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("D40").Value = "`'0.996"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").Value = "`'2.03"
$ws.Range("E42").Value = "  +8.39%  "
$ws.Range("E43").Value = "  +3.40%  "
$ws.Range("D44").Value = "`'134.57"
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("D45").Value = "`'17.24"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").Value = "`'3.89"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "`'22.44"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "`'2.21"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "2.209.62"
$ws.Range("E50").Value = "  +3.06%  "
$ws.Range("D51").Value = "3.796.32"
$ws.Range("E51").Value = "  +5.93%  "
